# Applies the ESG keyword list update:
#  - Column B: inserts "Quilombola" (social) in alphabetical order, shifting later rows down by 1
#  - Column C: inserts "Denuncia", "Golpe", "Ponzi", "Piramide financeira" (governance) in alphabetical
#    order, shifting later rows down (net +4 rows)
#  - Updates the hidden _xlnm._FilterDatabase defined name to reflect the new range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: shift rows 61-70 down (new word "Quilombola" inserted at B61) ---
$ws.Range('B61').Value = 'Quilombola'
$ws.Range('B62').Value = 'Relações trabalhistas'
$ws.Range('B63').Value = 'Responsabilidade social'
$ws.Range('B64').Value = 'Saúde'
$ws.Range('B65').Value = 'Sindicato'
$ws.Range('B66').Value = 'Trabalho escravo'
$ws.Range('B67').Value = 'Trabalho infantil'
$ws.Range('B68').Value = 'Tráfico'
$ws.Range('B69').Value = 'Transgêneros'
$ws.Range('B70').Value = 'Violência'

# --- Column C: shift rows 42-113 down (new words inserted at C42, C60, C89, C91) ---
$ws.Range('C42').Value = 'Denúncia'
$ws.Range('C43').Value = 'Desleal'
$ws.Range('C44').Value = 'Desvio conduta'
$ws.Range('C45').Value = 'Desvio dinheiro'
$ws.Range('C46').Value = 'Dividendos'
$ws.Range('C47').Value = 'Efeito negativo'
$ws.Range('C48').Value = 'Erro balanço'
$ws.Range('C49').Value = 'Erros contábeis'
$ws.Range('C50').Value = 'Escândalo'
$ws.Range('C51').Value = 'Ética'
$ws.Range('C52').Value = 'Evasão divisas'
$ws.Range('C53').Value = 'Evasão fiscal'
$ws.Range('C54').Value = 'Extrajudicial'
$ws.Range('C55').Value = 'Falência'
$ws.Range('C56').Value = 'Fraude'
$ws.Range('C57').Value = 'Furto'
$ws.Range('C58').Value = 'Gerenciamento Crise'
$ws.Range('C59').Value = 'Gestão riscos'
$ws.Range('C60').Value = 'Golpe'
$ws.Range('C61').Value = 'Governança'
$ws.Range('C62').Value = 'Honestidade'
$ws.Range('C63').Value = 'Ilegal'
$ws.Range('C64').Value = 'Ilícito'
$ws.Range('C65').Value = 'Incidente cibernético'
$ws.Range('C66').Value = 'Informação privilegiada'
$ws.Range('C67').Value = 'Inquérito'
$ws.Range('C68').Value = 'Insider trading'
$ws.Range('C69').Value = 'Insolvência'
$ws.Range('C70').Value = 'Inspeção'
$ws.Range('C71').Value = 'Investigação'
$ws.Range('C72').Value = 'Irregularidades'
$ws.Range('C73').Value = 'Justiça'
$ws.Range('C74').Value = 'Laranja'
$ws.Range('C75').Value = 'Lavagem dinheiro'
$ws.Range('C76').Value = 'LGPD'
$ws.Range('C77').Value = 'Lobby'
$ws.Range('C78').Value = 'Lobbysta'
$ws.Range('C79').Value = 'Lobista'
$ws.Range('C80').Value = 'Ministério Público'
$ws.Range('C81').Value = 'Minoritários'
$ws.Range('C82').Value = 'MPF'
$ws.Range('C83').Value = 'Multa'
$ws.Range('C84').Value = 'Negligência'
$ws.Range('C85').Value = 'Ocultação'
$ws.Range('C86').Value = 'Partes relacionadas'
$ws.Range('C87').Value = 'Partido político'
$ws.Range('C88').Value = 'Pessoas politicamente expostas'
$ws.Range('C89').Value = 'Pirâmide financeira'
$ws.Range('C90').Value = 'Polêmica'
$ws.Range('C91').Value = 'Ponzi'
$ws.Range('C92').Value = 'Prejuízo'
$ws.Range('C93').Value = 'Preso'
$ws.Range('C94').Value = 'Prisão'
$ws.Range('C95').Value = 'Privacidade'
$ws.Range('C96').Value = 'Propina'
$ws.Range('C97').Value = 'Punição'
$ws.Range('C98').Value = 'Recuperação judicial'
$ws.Range('C99').Value = 'Remuneração'
$ws.Range('C100').Value = 'Réu'
$ws.Range('C101').Value = 'Rombo contábil'
$ws.Range('C102').Value = 'Roubo'
$ws.Range('C103').Value = 'Sabotagem'
$ws.Range('C104').Value = 'Segurança cibernética'
$ws.Range('C105').Value = 'Segurança Dados'
$ws.Range('C106').Value = 'Skimming'
$ws.Range('C107').Value = 'Sonegação'
$ws.Range('C108').Value = 'Stakeholders'
$ws.Range('C109').Value = 'Suborno'
$ws.Range('C110').Value = 'Termo Ajustamento Conduta'
$ws.Range('C111').Value = 'Transparência'
$ws.Range('C112').Value = 'Vazamento dados'
$ws.Range('C113').Value = 'Violações'

# --- Update the hidden AutoFilter defined name range to match the edited extent ---
$fdb = $wb.Names.Item('_xlnm._FilterDatabase')
$fdb.RefersTo = '=Planilha1!$C$1:$C$67'

# --- Restore the view/selection state (scrolled down a bit further after the new rows) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$ws.Range('C114').Select() | Out-Null
